# COREESG_holdings.xlsx - "Add files via upload" re-upload with refreshed
# model weights/percent-changes and a one-day-later "as of" date in the
# confidentiality footer.
#
# The sheet ships protected (password "D382"), so it must be unprotected
# before the cell values can be written, and re-protected with the same
# password afterwards to restore the original protected state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect("D382")

# Weight (D) / Percent Change (E) columns, rows 2-7 (NULG, NULV, NUMG, NUMV,
# NUSC, Total)
$ws.Range("D2").Value = 0.2542670607410571
$ws.Range("E2").Value = -0.004163997437540035

$ws.Range("D3").Value = 0.4870337071313253
$ws.Range("E3").Value = 0

$ws.Range("D4").Value = 0.101876457255036
$ws.Range("E4").Value = -0.003883854262992426

$ws.Range("D5").Value = 0.09959390069531461
$ws.Range("E5").Value = 0.00330669605952072

$ws.Range("D6").Value = 0.05722887417726688
$ws.Range("E6").Value = 0.002663115845539243

$ws.Range("D7").Value = 0.9999999999999999
$ws.Range("E7").Value = -0.00097270682155981

# Confidentiality footer: bump the "Model holdings provided as of" date by
# one day (2021-04-26 -> 2021-04-27).
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-27 for illustrative purposes only and are subject to change."

$ws.Protect("D382")
